$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = 5985
$ws.Range("B16").Value = "bert_classifier_5985.pth"
$ws.Range("C16").Value = "5985_accuracy_scores.docx"
$ws.Range("D16").Value = 1773
$ws.Range("E16").Value = 846
$ws.Range("F16").Value = 533
$ws.Range("G16").Value = 394
$ws.Range("H16").Value = "bert-large-cased-whole-word-masking"
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 0.00002
$ws.Range("K16").Value = 16
$ws.Range("L16").Value = 512
$ws.Range("M16").Value = 0.952112676056338
$ws.Range("N16").Value = 10

# Row 17
$ws.Range("A17").Value = 5964
$ws.Range("B17").Value = "bert_classifier_5964.pth"
$ws.Range("C17").Value = "5964_accuracy_scores.docx"
$ws.Range("D17").Value = 1773
$ws.Range("E17").Value = 846
$ws.Range("F17").Value = 533
$ws.Range("G17").Value = 394
$ws.Range("H17").Value = "bert-base-uncased"
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 0.00005
$ws.Range("K17").Value = 32
$ws.Range("L17").Value = 256
$ws.Range("M17").Value = 0.9323943661971831
$ws.Range("N17").Value = 10

# Row 18
$ws.Range("A18").Value = 5537
$ws.Range("B18").Value = "bert_classifier_5537.pth"
$ws.Range("C18").Value = "5537_accuracy_scores.docx"
$ws.Range("D18").Value = 1773
$ws.Range("E18").Value = 846
$ws.Range("F18").Value = 533
$ws.Range("G18").Value = 394
$ws.Range("H18").Value = "bert-base-uncased"
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 0.0003
$ws.Range("K18").Value = 16
$ws.Range("L18").Value = 512
$ws.Range("M18").Value = 0.5267605633802817
$ws.Range("N18").Value = 10

# Row 19
$ws.Range("A19").Value = 4223
$ws.Range("B19").Value = "bert_classifier_4223.pth"
$ws.Range("C19").Value = "4223_accuracy_scores.docx"
$ws.Range("D19").Value = 6412
$ws.Range("E19").Value = 2690
$ws.Range("F19").Value = 2140
$ws.Range("G19").Value = 1582
$ws.Range("H19").Value = "bert-base-cased"
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 0.0003
$ws.Range("K19").Value = 8
$ws.Range("L19").Value = 256
$ws.Range("M19").Value = 0.436477007014809
$ws.Range("N19").Value = 25

# Row 20
$ws.Range("A20").Value = 3142
$ws.Range("B20").Value = "bert_classifier_3142.pth"
$ws.Range("C20").Value = "3142_accuracy_scores.docx"
$ws.Range("D20").Value = 5889
$ws.Range("E20").Value = 2499
$ws.Range("F20").Value = 1963
$ws.Range("G20").Value = 1427
$ws.Range("H20").Value = "bert-large-uncased"
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 0.0001
$ws.Range("K20").Value = 16
$ws.Range("L20").Value = 256
$ws.Range("M20").Value = 0.4320882852292021
$ws.Range("N20").Value = 25

# Row 21
$ws.Range("A21").Value = 805
$ws.Range("B21").Value = "bert_classifier_805.pth"
$ws.Range("C21").Value = "805_accuracy_scores.docx"
$ws.Range("D21").Value = 6412
$ws.Range("E21").Value = 2690
$ws.Range("F21").Value = 2140
$ws.Range("G21").Value = 1582
$ws.Range("H21").Value = "bert-base-uncased"
$ws.Range("I21").Value = 4
$ws.Range("J21").Value = 0.00003
$ws.Range("K21").Value = 16
$ws.Range("L21").Value = 512
$ws.Range("M21").Value = 0.9220576773187841
$ws.Range("N21").Value = 25

# Row 22
$ws.Range("A22").Value = 2344
$ws.Range("B22").Value = "bert_classifier_2344.pth"
$ws.Range("C22").Value = "2344_accuracy_scores.docx"
$ws.Range("D22").Value = 6412
$ws.Range("E22").Value = 2690
$ws.Range("F22").Value = 2140
$ws.Range("G22").Value = 1582
$ws.Range("H22").Value = "bert-base-cased"
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 0.0003
$ws.Range("K22").Value = 64
$ws.Range("L22").Value = 256
$ws.Range("M22").Value = 0.436477007014809
$ws.Range("N22").Value = 25

# Row 23
$ws.Range("A23").Value = 7918
$ws.Range("B23").Value = "bert_classifier_7918.pth"
$ws.Range("C23").Value = "7918_accuracy_scores.docx"
$ws.Range("D23").Value = 1773
$ws.Range("E23").Value = 846
$ws.Range("F23").Value = 533
$ws.Range("G23").Value = 394
$ws.Range("H23").Value = "bert-large-cased-whole-word-masking"
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 0.00005
$ws.Range("K23").Value = 16
$ws.Range("L23").Value = 512
$ws.Range("M23").Value = 0.5267605633802817
$ws.Range("N23").Value = 10

# Row 24
$ws.Range("A24").Value = 5909
$ws.Range("B24").Value = "bert_classifier_5909.pth"
$ws.Range("C24").Value = "5909_accuracy_scores.docx"
$ws.Range("D24").Value = 1773
$ws.Range("E24").Value = 846
$ws.Range("F24").Value = 533
$ws.Range("G24").Value = 394
$ws.Range("H24").Value = "bert-large-uncased"
$ws.Range("I24").Value = 6
$ws.Range("J24").Value = 0.0001
$ws.Range("K24").Value = 8
$ws.Range("L24").Value = 256
$ws.Range("M24").Value = 0.5267605633802817
$ws.Range("N24").Value = 10

# Row 25
$ws.Range("A25").Value = 2290
$ws.Range("B25").Value = "bert_classifier_2290.pth"
$ws.Range("C25").Value = "2290_accuracy_scores.docx"
$ws.Range("D25").Value = 1773
$ws.Range("E25").Value = 846
$ws.Range("F25").Value = 533
$ws.Range("G25").Value = 394
$ws.Range("H25").Value = "bert-large-cased"
$ws.Range("I25").Value = 4
$ws.Range("J25").Value = 0.00002
$ws.Range("K25").Value = 16
$ws.Range("L25").Value = 512
$ws.Range("M25").Value = 0.9352112676056338
$ws.Range("N25").Value = 10

# Row 26
$ws.Range("A26").Value = 8978
$ws.Range("B26").Value = "bert_classifier_8978.pth"
$ws.Range("C26").Value = "8978_accuracy_scores.docx"
$ws.Range("D26").Value = 1773
$ws.Range("E26").Value = 846
$ws.Range("F26").Value = 533
$ws.Range("G26").Value = 394
$ws.Range("H26").Value = "bert-large-uncased"
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 0.00005
$ws.Range("K26").Value = 8
$ws.Range("L26").Value = 256
$ws.Range("M26").Value = 0.5267605633802817
$ws.Range("N26").Value = 10

# Row 27
$ws.Range("A27").Value = 3238
$ws.Range("B27").Value = "bert_classifier_3238.pth"
$ws.Range("C27").Value = "3238_accuracy_scores.docx"
$ws.Range("D27").Value = 1773
$ws.Range("E27").Value = 846
$ws.Range("F27").Value = 533
$ws.Range("G27").Value = 394
$ws.Range("H27").Value = "bert-base-uncased"
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 0.00002
$ws.Range("K27").Value = 32
$ws.Range("L27").Value = 512
$ws.Range("M27").Value = 0.923943661971831
$ws.Range("N27").Value = 10

